$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:D2").Value = 15.45
$ws.Range("B3:D3").Value = 932.8199999999999
$ws.Range("B4:C4").Value = 798.1
$ws.Range("D4").Value = 757.15
